# Apply the betexplorer Serie A 2023-2024 refresh:
#  - six pairs of adjacent fixtures had their match details (columns F:V)
#    swapped (their scrape order changed upstream; the row's
#    Indice/pais/torneio/temporada/data_partida in A..E stay put)
#  - four new fixtures (rows 162-165) were appended at the bottom
#  - the sheet <dimension> grows from A1:V161 to A1:V165 (handled
#    automatically by the engine once the new cells are written)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- swap the "home..url" data (columns F..V) between two rows, keeping
#     Indice/pais/torneio/temporada/data_partida (A..E) where they are.
#     (Positional params only -- this host's PowerShell-alike chokes on
#     named/-switch style arguments to user functions.) ---
function Swap-MatchDetails {
    param($row1, $row2)

    for ($col = 6; $col -le 22; $col++) {
        $val1 = $ws.Cells.Item($row1, $col).Value2
        $val2 = $ws.Cells.Item($row2, $col).Value2
        $ws.Cells.Item($row1, $col).Value = $val2
        $ws.Cells.Item($row2, $col).Value = $val1
    }
}

Swap-MatchDetails 6 7
Swap-MatchDetails 26 27
Swap-MatchDetails 118 119
Swap-MatchDetails 126 127
Swap-MatchDetails 137 138
Swap-MatchDetails 157 158

# --- append the four newly scraped fixtures as rows 162..165 ---
function Set-MatchRow {
    param(
        $row, $indice, $dataPartida,
        $home, $homeGols, $away, $awayGols,
        $homeOpenOdds, $homeOpenDh, $homeCloseOdds, $homeCloseDh,
        $drawOpenOdds, $drawOpenDh, $drawCloseOdds, $drawCloseDh,
        $awayOpenOdds, $awayOpenDh, $awayCloseOdds, $awayCloseDh,
        $url
    )

    $ws.Cells.Item($row, 1).Value = $indice
    $ws.Cells.Item($row, 2).Value = "italy"
    $ws.Cells.Item($row, 3).Value = "serie-a"
    $ws.Cells.Item($row, 4).Value = "2023-2024"
    $ws.Cells.Item($row, 5).Value = $dataPartida
    $ws.Cells.Item($row, 6).Value = $home
    $ws.Cells.Item($row, 7).Value = $homeGols
    $ws.Cells.Item($row, 8).Value = $away
    $ws.Cells.Item($row, 9).Value = $awayGols
    $ws.Cells.Item($row, 10).Value = $homeOpenOdds
    $ws.Cells.Item($row, 11).Value = $homeOpenDh
    $ws.Cells.Item($row, 12).Value = $homeCloseOdds
    $ws.Cells.Item($row, 13).Value = $homeCloseDh
    $ws.Cells.Item($row, 14).Value = $drawOpenOdds
    $ws.Cells.Item($row, 15).Value = $drawOpenDh
    $ws.Cells.Item($row, 16).Value = $drawCloseOdds
    $ws.Cells.Item($row, 17).Value = $drawCloseDh
    $ws.Cells.Item($row, 18).Value = $awayOpenOdds
    $ws.Cells.Item($row, 19).Value = $awayOpenDh
    $ws.Cells.Item($row, 20).Value = $awayCloseOdds
    $ws.Cells.Item($row, 21).Value = $awayCloseDh
    $ws.Cells.Item($row, 22).Value = $url

    # Match the formatting of the existing data rows: bold/bordered/centered
    # index in column A, date-time number format in column E. Copy it from
    # the row just above (already scraped data) rather than re-declaring a
    # style so the workbook doesn't grow a near-duplicate style entry.
    $prevRow = $row - 1
    $ws.Cells.Item($prevRow, 1).Copy() | Out-Null
    $ws.Cells.Item($row, 1).PasteSpecial(-4122) | Out-Null
    $ws.Cells.Item($prevRow, 5).Copy() | Out-Null
    $ws.Cells.Item($row, 5).PasteSpecial(-4122) | Out-Null
}

Set-MatchRow 162 161 45282.77083333334 `
    "Empoli" 0 "Lazio" 2 `
    3.36 "10/12/2023 10:02" 4.61 "22/12/2023 18:29" `
    3.49 "10/12/2023 10:02" 3.52 "22/12/2023 18:29" `
    2.08 "10/12/2023 10:02" 1.93 "22/12/2023 17:43" `
    "https://www.betexplorer.com/football/italy/serie-a/empoli-lazio/ImZ8Us9N/"

Set-MatchRow 163 162 45282.77083333334 `
    "Sassuolo" 1 "Genoa" 2 `
    2 "10/12/2023 10:02" 2.59 "22/12/2023 18:29" `
    3.54 "10/12/2023 10:02" 3.1 "22/12/2023 18:28" `
    3.52 "10/12/2023 10:02" 3.14 "22/12/2023 18:29" `
    "https://www.betexplorer.com/football/italy/serie-a/sassuolo-genoa/jyrhYP9b/"

Set-MatchRow 164 163 45282.86458333334 `
    "Salernitana" 2 "AC Milan" 2 `
    5.3 "10/12/2023 10:02" 6.64 "22/12/2023 20:44" `
    4.1 "10/12/2023 10:02" 4.73 "22/12/2023 20:44" `
    1.56 "10/12/2023 10:02" 1.49 "22/12/2023 20:43" `
    "https://www.betexplorer.com/football/italy/serie-a/salernitana-ac-milan/CpqlZ5fh/"

Set-MatchRow 165 164 45282.86458333334 `
    "Monza" 0 "Fiorentina" 1 `
    2.77 "10/12/2023 10:02" 2.89 "22/12/2023 20:40" `
    3.32 "10/12/2023 10:02" 3.34 "22/12/2023 20:44" `
    2.49 "10/12/2023 10:02" 2.63 "22/12/2023 20:37" `
    "https://www.betexplorer.com/football/italy/serie-a/monza-fiorentina/4KutySPu/"
